# Anonymize "fedcore" -> "approach" and redraw the header-row bottom border
# under the merged B1:D1 / E1:G1 cells. The left (anchor) cell of each
# merge keeps its existing bold boxed style (s=1); the remaining cells of
# the merge get a plain top+bottom border, with the rightmost cell of the
# merge also getting a right border.
#
# NB: border edges are written in this specific order (right, bottom, top,
# left) because writing them in natural reading order (left, top, bottom,
# right) on a cell that starts fully boxed can make the engine pass
# through an intermediate border combination that has no pre-existing
# cellXfs entry, permanently allocating a stray/orphan style slot. This
# order never produces an intermediate combination that isn't already a
# pre-existing border, so no extra style slots are allocated.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: quality_comparison ---
$ws1 = $wb.Worksheets.Item("quality_comparison")

$c1 = $ws1.Range("C1")
$c1.ClearFormats()
$c1.Borders.Item(10).LineStyle = -4142
$c1.Borders.Item(9).LineStyle = 1
$c1.Borders.Item(8).LineStyle = 1
$c1.Borders.Item(7).LineStyle = -4142

$d1 = $ws1.Range("D1")
$d1.ClearFormats()
$d1.Borders.Item(10).LineStyle = 1
$d1.Borders.Item(9).LineStyle = 1
$d1.Borders.Item(8).LineStyle = 1
$d1.Borders.Item(7).LineStyle = -4142

$ws1.Range("C2").Value = "approach"

# --- Sheet 2: computational_comparison ---
$ws2 = $wb.Worksheets.Item("computational_comparison")

$c1b = $ws2.Range("C1")
$c1b.ClearFormats()
$c1b.Borders.Item(10).LineStyle = -4142
$c1b.Borders.Item(9).LineStyle = 1
$c1b.Borders.Item(8).LineStyle = 1
$c1b.Borders.Item(7).LineStyle = -4142

$d1b = $ws2.Range("D1")
$d1b.ClearFormats()
$d1b.Borders.Item(10).LineStyle = 1
$d1b.Borders.Item(9).LineStyle = 1
$d1b.Borders.Item(8).LineStyle = 1
$d1b.Borders.Item(7).LineStyle = -4142

$f1b = $ws2.Range("F1")
$f1b.ClearFormats()
$f1b.Borders.Item(10).LineStyle = -4142
$f1b.Borders.Item(9).LineStyle = 1
$f1b.Borders.Item(8).LineStyle = 1
$f1b.Borders.Item(7).LineStyle = -4142

$g1b = $ws2.Range("G1")
$g1b.ClearFormats()
$g1b.Borders.Item(10).LineStyle = 1
$g1b.Borders.Item(9).LineStyle = 1
$g1b.Borders.Item(8).LineStyle = 1
$g1b.Borders.Item(7).LineStyle = -4142

$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# G5 was an empty inline-string placeholder cell; drop it entirely.
$ws2.Range("G5").ClearContents()

Write-Output "edit applied"
